# Insert a new data row before row 109 (shifting existing rows 109-121 down to 110-122)
# and populate the newly-inserted row 109 with the new record described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 109..121 down to 110..122, creating a blank row 109.
$ws.Rows.Item(109).Insert()

# Fill the new row 109 with the new record's data.
$ws.Range("A109").Value = 1
$ws.Range("B109").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C109").Value = "Arica y Parinacota"
$ws.Range("D109").Value = 44951
$ws.Range("E109").Value = 15
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100102
$ws.Range("H109").Value = "Cítricos"
$ws.Range("I109").Value = 100102005
$ws.Range("J109").Value = "Naranja"
$ws.Range("K109").Value = "Valencia"
$ws.Range("L109").Value = "Segunda"
$ws.Range("M109").Value = 650
$ws.Range("N109").Value = 950
$ws.Range("O109").Value = 1000
$ws.Range("P109").Value = 977
$ws.Range("Q109").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R109").Value = "Región de Coquimbo"
$ws.Range("S109").Value = 977
$ws.Range("T109").Value = 1
